# Auto-generated Excel COM-interop edit script
# Applies numeric cell updates (and a few cell additions/removals) across the 8 class worksheets
# per the authoritative diff of Sheets/Halicarnassus_Profits.xlsx

$wb = $excel.ActiveWorkbook

# ---- Worksheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value2 = 2157.976  # was 1743.0444
$ws.Range("I15").Value2 = 2157.976  # was 1743.0444
$ws.Range("K15").Value2 = 6473.928  # was 5229.1332
$ws.Range("M15").Value2 = -6304.928  # was -5060.1332
$ws.Range("H113").Value2 = 2726.75  # was 3235.3333
$ws.Range("I113").Value2 = 3100.5  # was 5000
$ws.Range("K113").Value2 = 3100.5  # was 5000
$ws.Range("M113").Value2 = 153.5  # was -1746
$ws.Range("H125").Value2 = 2935.45  # was 3208.3684
$ws.Range("I125").Value2 = 3047.8462  # was 3336.3076
$ws.Range("J125").Value2 = 2726.7144  # was 2931.1667
$ws.Range("K125").Value2 = 27430.6158  # was 30026.7684
$ws.Range("L125").Value2 = 24540.4296  # was 26380.5003
$ws.Range("M125").Value2 = -24970.6158  # was -27566.7684
$ws.Range("N125").Value2 = -29460.4296  # was -31300.5003
$ws.Range("H138").Value2 = 2324.2222  # was 2131.65
$ws.Range("I138").Value2 = 808.5  # was 757.25
$ws.Range("K138").Value2 = 2425.5  # was 2271.75
$ws.Range("M138").Value2 = 2714.5  # was 2868.25

# ---- Worksheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value2 = 1277.75  # was 1217.1538
$ws.Range("I61").Value2 = 1164.8889  # was 1097.4
$ws.Range("K61").Value2 = 1164.8889  # was 1097.4
$ws.Range("M61").Value2 = -952.8888999999999  # was -885.4000000000001
$ws.Range("H74").Value2 = 3480.1052  # was 3067.3333
$ws.Range("I74").Value2 = 3223.7222  # was 2815.95
$ws.Range("K74").Value2 = 3223.7222  # was 2815.95
$ws.Range("M74").Value2 = -2349.7222  # was -1941.95
$ws.Range("H77").Value2 = 3480.1052  # was 3067.3333
$ws.Range("I77").Value2 = 3223.7222  # was 2815.95
$ws.Range("K77").Value2 = 16118.611  # was 14079.75
$ws.Range("M77").Value2 = -11750.611  # was -9711.75
$ws.Range("H132").Value2 = 2384.3684  # was 2462
$ws.Range("I132").Value2 = 2295.75  # was 2414.7273
$ws.Range("K132").Value2 = 6887.25  # was 7244.1819
$ws.Range("M132").Value2 = -4357.25  # was -4714.1819
$ws.Range("H134").Value2 = 50000  # was 0
$ws.Range("J134").Value2 = 50000  # was 0
$ws.Range("L134").Value2 = 50000  # was 0
$ws.Range("N134").Value2 = -60140  # new cell
$ws.Range("H136").Value2 = 1277.75  # was 1217.1538
$ws.Range("I136").Value2 = 1164.8889  # was 1097.4
$ws.Range("K136").Value2 = 3494.6667  # was 3292.2
$ws.Range("M136").Value2 = -944.6666999999998  # was -742.2000000000003

# ---- Worksheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H44").Value2 = 30050  # was 30049.334
$ws.Range("J44").Value2 = 30050  # was 30049.334
$ws.Range("L44").Value2 = 30050  # was 30049.334
$ws.Range("N44").Value2 = -31044  # was -31043.334
$ws.Range("H99").Value2 = 2120  # was 2175
$ws.Range("I99").Value2 = 2105  # was 2200
$ws.Range("K99").Value2 = 2105  # was 2200
$ws.Range("M99").Value2 = -607  # was -702
$ws.Range("H107").Value2 = 4333.3335  # was 375.16666
$ws.Range("I107").Value2 = 4000  # was 371.2
$ws.Range("J107").Value2 = 5000  # was 395
$ws.Range("K107").Value2 = 4000  # was 371.2
$ws.Range("L107").Value2 = 5000  # was 395
$ws.Range("M107").Value2 = -2080  # was 1548.8
$ws.Range("N107").Value2 = -8840  # was -4235
$ws.Range("H134").Value2 = 4439.2  # was 4254
$ws.Range("I134").Value2 = 1510.2222  # was 1459.6842
$ws.Range("K134").Value2 = 4530.6666  # was 4379.0526
$ws.Range("M134").Value2 = -1995.6666  # was -1844.0526

# ---- Worksheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value2 = 1039.6  # was 1011.4375
$ws.Range("I22").Value2 = 935.5714  # was 892.25
$ws.Range("K22").Value2 = 935.5714  # was 892.25
$ws.Range("M22").Value2 = -585.5714  # was -542.25
$ws.Range("H44").Value2 = 19996  # was 19997.5
$ws.Range("J44").Value2 = 0  # was 19999
$ws.Range("L44").Value2 = 0  # was 19999
$ws.Range("N44").ClearContents()  # was -20883
$ws.Range("H45").Value2 = 0  # was 34999
$ws.Range("J45").Value2 = 0  # was 34999
$ws.Range("L45").Value2 = 0  # was 34999
$ws.Range("N45").ClearContents()  # was -36185
$ws.Range("H47").Value2 = 20000  # was 31666
$ws.Range("I47").Value2 = 20000  # was 25000
$ws.Range("J47").Value2 = 0  # was 34999
$ws.Range("K47").Value2 = 20000  # was 25000
$ws.Range("L47").Value2 = 0  # was 34999
$ws.Range("M47").Value2 = -19434  # was -24434
$ws.Range("N47").ClearContents()  # was -36131
$ws.Range("H55").Value2 = 5059.75  # was 9891.286
$ws.Range("I55").Value2 = 5386  # was 8999
$ws.Range("J55").Value2 = 4081  # was 10248.2
$ws.Range("K55").Value2 = 5386  # was 8999
$ws.Range("L55").Value2 = 4081  # was 10248.2
$ws.Range("M55").Value2 = -5071  # was -8684
$ws.Range("N55").Value2 = -4711  # was -10878.2
$ws.Range("H58").Value2 = 3273.889  # was 3116.2
$ws.Range("I58").Value2 = 2130  # was 2079.0588
$ws.Range("K58").Value2 = 2130  # was 2079.0588
$ws.Range("M58").Value2 = -1927  # was -1876.0588
$ws.Range("H99").Value2 = 2290  # was 3500
$ws.Range("I99").Value2 = 1975  # was 2875
$ws.Range("J99").Value2 = 3550  # was 6000
$ws.Range("K99").Value2 = 1975  # was 2875
$ws.Range("L99").Value2 = 3550  # was 6000
$ws.Range("M99").Value2 = -477  # was -1377
$ws.Range("N99").Value2 = -6546  # was -8996
$ws.Range("H105").Value2 = 1213  # was 1510.6666
$ws.Range("J105").Value2 = 2900  # was 2949.5
$ws.Range("L105").Value2 = 2900  # was 2949.5
$ws.Range("N105").Value2 = -6394  # was -6443.5
$ws.Range("H126").Value2 = 2290  # was 3500
$ws.Range("I126").Value2 = 1975  # was 2875
$ws.Range("J126").Value2 = 3550  # was 6000
$ws.Range("K126").Value2 = 5925  # was 8625
$ws.Range("L126").Value2 = 10650  # was 18000
$ws.Range("M126").Value2 = -3455  # was -6155
$ws.Range("N126").Value2 = -15590  # was -22940
$ws.Range("H132").Value2 = 2712.3416  # was 2966.4167
$ws.Range("I132").Value2 = 2523.7058  # was 2806.5862
$ws.Range("K132").Value2 = 7571.117400000001  # was 8419.758600000001
$ws.Range("M132").Value2 = -5041.117400000001  # was -5889.758600000001
$ws.Range("H134").Value2 = 2009.9524  # was 2147.5789
$ws.Range("I134").Value2 = 1274.1578  # was 1341.4117
$ws.Range("K134").Value2 = 3822.4734  # was 4024.2351
$ws.Range("M134").Value2 = -1287.4734  # was -1489.2351
$ws.Range("H136").Value2 = 3273.889  # was 3116.2
$ws.Range("I136").Value2 = 2130  # was 2079.0588
$ws.Range("K136").Value2 = 6390  # was 6237.176399999999
$ws.Range("M136").Value2 = -3840  # was -3687.176399999999

# ---- Worksheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H11").Value2 = 81066.336  # was 86855
$ws.Range("I11").Value2 = 100566.25  # was 109706.37
$ws.Range("K11").Value2 = 301698.75  # was 329119.11
$ws.Range("M11").Value2 = -301558.75  # was -328979.11
$ws.Range("H34").Value2 = 1602.7059  # was 1708.1111
$ws.Range("J34").Value2 = 3535.1428  # was 3530.75
$ws.Range("L34").Value2 = 10605.4284  # was 10592.25
$ws.Range("N34").Value2 = -10773.4284  # was -10760.25
$ws.Range("H39").Value2 = 8612.25  # was 8389.799999999999
$ws.Range("J39").Value2 = 8612.25  # was 8389.799999999999
$ws.Range("L39").Value2 = 25836.75  # was 25169.4
$ws.Range("N39").Value2 = -26424.75  # was -25757.4

# ---- Worksheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H48").Value2 = 0  # was 10000
$ws.Range("J48").Value2 = 0  # was 10000
$ws.Range("L48").Value2 = 0  # was 10000
$ws.Range("N48").ClearContents()  # was -10970
$ws.Range("H132").Value2 = 35415.78  # was 34418.184
$ws.Range("I132").Value2 = 44025.04  # was 42427.73
$ws.Range("K132").Value2 = 132075.12  # was 127283.19
$ws.Range("M132").Value2 = -129545.12  # was -124753.19
$ws.Range("H140").Value2 = 141443.5  # was 128799.5
$ws.Range("J140").Value2 = 95359  # was 98620
$ws.Range("L140").Value2 = 95359  # was 98620
$ws.Range("N140").Value2 = -105719  # was -108980

# ---- Worksheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("J7").Value2 = 3000  # was 0
$ws.Range("L7").Value2 = 3000  # was 0
$ws.Range("N7").Value2 = -3224  # new cell
$ws.Range("H46").Value2 = 2039.2858  # was 2659
$ws.Range("I46").Value2 = 493.33334  # was 500
$ws.Range("K46").Value2 = 493.33334  # was 500
$ws.Range("M46").Value2 = -305.33334  # was -312
$ws.Range("H82").Value2 = 1800  # was 2083.3333
$ws.Range("I82").Value2 = 1800  # was 2083.3333
$ws.Range("K82").Value2 = 1800  # was 2083.3333
$ws.Range("M82").Value2 = -1439  # was -1722.3333
$ws.Range("H85").Value2 = 1800  # was 2083.3333
$ws.Range("I85").Value2 = 1800  # was 2083.3333
$ws.Range("K85").Value2 = 1800  # was 2083.3333
$ws.Range("M85").Value2 = -552  # was -835.3332999999998
$ws.Range("J126").Value2 = 3000  # was 0
$ws.Range("L126").Value2 = 9000  # was 0
$ws.Range("N126").Value2 = -13940  # new cell
$ws.Range("H136").Value2 = 3349.2  # was 2961.8333
$ws.Range("I136").Value2 = 2064.889  # was 1744.9231
$ws.Range("K136").Value2 = 6194.667  # was 5234.7693
$ws.Range("M136").Value2 = -3644.667  # was -2684.7693

# ---- Worksheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H52").Value2 = 20760.5  # was 24013.666
$ws.Range("I52").Value2 = 12680.667  # was 13021
$ws.Range("J52").Value2 = 45000  # was 45999
$ws.Range("K52").Value2 = 12680.667  # was 13021
$ws.Range("L52").Value2 = 45000  # was 45999
$ws.Range("M52").Value2 = -12454.667  # was -12795
$ws.Range("N52").Value2 = -45452  # was -46451
$ws.Range("H122").Value2 = 2049.5334  # was 2105.2068
$ws.Range("I122").Value2 = 1625.4584  # was 1677.2174
$ws.Range("K122").Value2 = 4876.3752  # was 5031.6522
$ws.Range("M122").Value2 = -2426.3752  # was -2581.6522
$ws.Range("H126").Value2 = 3650.077  # was 3296.7334
$ws.Range("J126").Value2 = 5583.3335  # was 3750
$ws.Range("L126").Value2 = 16750.0005  # was 11250
$ws.Range("N126").Value2 = -21690.0005  # was -16190
$ws.Range("H132").Value2 = 3471.1428  # was 3472
$ws.Range("J132").Value2 = 3499.5  # was 3484.2
$ws.Range("L132").Value2 = 10498.5  # was 10452.6
$ws.Range("N132").Value2 = -15558.5  # was -15512.6
$ws.Range("H136").Value2 = 2730.5  # was 2594.639
$ws.Range("I136").Value2 = 2109.48  # was 1974.3334
$ws.Range("K136").Value2 = 6328.440000000001  # was 5923.0002
$ws.Range("M136").Value2 = -3778.440000000001  # was -3373.0002
